# Pruebas funcionales TC 003 Alcosto
# Add a new worksheet "Data003" with search term data, make it the active sheet,
# and clear the previous selection/active-cell state on "Data002".

$wb = $excel.ActiveWorkbook

# Clear tabSelected/activeCell state on Data002 (second sheet) - select full range A1:B5
$ws2 = $wb.Worksheets.Item("Data002")
$ws2.Range("A1:B5").Select() | Out-Null

# Add new sheet "Data003" right after Data002 (becomes the last / active tab)
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Data003"

# Populate new sheet with data
$ws3.Range("A1").Value = "busqueda"
$ws3.Range("A2").Value = "Portatil Lenovo"

# Auto-fit column A to the new content (mirrors the authoring session)
$ws3.Columns.Item(1).AutoFit() | Out-Null

# Select A2 on the new sheet, and make it the active sheet/tab
$ws3.Select() | Out-Null
$ws3.Range("A2").Select() | Out-Null
